$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price/percentage text columns keep their original text formatting
# (avoid Excel auto-converting numeric-looking strings into numbers)
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "30.136.15"
$ws.Range("E2").Value = "  +5.37%  "

# Row 3
$ws.Range("D3").Value = "1.921.61"
$ws.Range("E3").Value = "  +2.19%  "

# Row 4
$ws.Range("E4").Value = "  -1.09%  "

# Row 5
$ws.Range("D5").Value = "327.55"
$ws.Range("E5").Value = "  +3.44%  "

# Row 6
$ws.Range("E6").Value = "  -1.03%  "

# Row 7
$ws.Range("D7").Value = "0.5162"
$ws.Range("E7").Value = "  +1.25%  "

# Row 8
$ws.Range("D8").Value = "0.4012"
$ws.Range("E8").Value = "  +2.69%  "

# Row 9
$ws.Range("D9").Value = "0.08461"
$ws.Range("E9").Value = "  +0.48%  "

# Row 10
$ws.Range("D10").Value = "42.75"
$ws.Range("E10").Value = "  +1.96%  "

# Row 11
$ws.Range("D11").Value = "1.124"
$ws.Range("E11").Value = "  +1.70%  "

# Row 12
$ws.Range("D12").Value = "21.55"
$ws.Range("E12").Value = "  +5.43%  "

# Row 13
$ws.Range("D13").Value = "6.349"
$ws.Range("E13").Value = "  +1.87%  "

# Row 14
$ws.Range("D14").Value = "1.920.96"
$ws.Range("E14").Value = "  +2.41%  "

# Row 15
$ws.Range("D15").Value = "7.360"
$ws.Range("E15").Value = "  +1.46%  "

# Row 16
$ws.Range("D16").Value = "1.000"
$ws.Range("E16").Value = "  -1.22%  "

# Row 17
$ws.Range("D17").Value = "96.29"
$ws.Range("E17").Value = "  +5.33%  "

# Row 18
$ws.Range("D18").Value = "0.00001117"
$ws.Range("E18").Value = "  +1.03%  "

# Row 19
$ws.Range("D19").Value = "0.06731"
$ws.Range("E19").Value = "  -0.05%  "

# Row 20
$ws.Range("D20").Value = "18.12"
$ws.Range("E20").Value = "  +2.11%  "

# Row 21
$ws.Range("E21").Value = "  -0.94%  "

# Row 22
$ws.Range("D22").Value = "6.071"
$ws.Range("E22").Value = "  +2.22%  "

# Row 23
$ws.Range("D23").Value = "30.156.25"
$ws.Range("E23").Value = "  +5.32%  "

# Row 24
$ws.Range("D24").Value = "11.23"
$ws.Range("E24").Value = "  +1.07%  "

# Row 25
$ws.Range("D25").Value = "2.203"
$ws.Range("E25").Value = "  -1.77%  "

# Row 26
$ws.Range("D26").Value = "2.140.41"
$ws.Range("E26").Value = "  +2.44%  "

# Row 27
$ws.Range("D27").Value = "160.77"
$ws.Range("E27").Value = "  -0.78%  "

# Row 28
$ws.Range("D28").Value = "21.04"
$ws.Range("E28").Value = "  +1.71%  "

# Row 29
$ws.Range("D29").Value = "2.465"
$ws.Range("E29").Value = "  +4.54%  "

# Row 30
$ws.Range("D30").Value = "129.09"
$ws.Range("E30").Value = "  +2.27%  "

# Row 31
$ws.Range("D31").Value = "1.079"
$ws.Range("E31").Value = "  +3.55%  "

# Row 32
$ws.Range("E32").Value = "  +1.35%  "

# Row 33
$ws.Range("D33").Value = "6.083"
$ws.Range("E33").Value = "  +5.05%  "

# Row 34
$ws.Range("D34").Value = "3.671"
$ws.Range("E34").Value = "  +1.24%  "

# Row 35
$ws.Range("D35").Value = "0.02515"
$ws.Range("E35").Value = "  +2.02%  "

# Row 36
$ws.Range("D36").Value = "0.06610"
$ws.Range("E36").Value = "  +0.82%  "

# Row 37
$ws.Range("D37").Value = "0.2226"
$ws.Range("E37").Value = "  +2.77%  "

# Row 38
$ws.Range("D38").Value = "1.237"
$ws.Range("E38").Value = "  +3.37%  "

# Row 39
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").Value = "9.021"
$ws.Range("E39").Value = "  +1.99%  "

# Row 40
$ws.Range("B40").Value = "InternetComputer(DFINITY)"
$ws.Range("C40").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D40").Value = "5.208"
$ws.Range("E40").Value = "  +2.40%  "

# Row 41
$ws.Range("D41").Value = "0.6564"
$ws.Range("E41").Value = "  +2.02%  "

# Row 42
$ws.Range("D42").Value = "1.243"
$ws.Range("E42").Value = "  -1.10%  "

# Row 43
$ws.Range("D43").Value = "11.42"
$ws.Range("E43").Value = "  +2.58%  "

# Row 44
$ws.Range("D44").Value = "0.6146"

# Row 45
$ws.Range("E45").Value = "  +0.74%  "

# Row 46
$ws.Range("D46").Value = "3.767"
$ws.Range("E46").Value = "  +1.83%  "

# Row 47
$ws.Range("D47").Value = "2.059"
$ws.Range("E47").Value = "  +2.38%  "

# Row 48
$ws.Range("D48").Value = "126.00"
$ws.Range("E48").Value = "  +3.06%  "

# Row 49
$ws.Range("D49").Value = "1.244"
$ws.Range("E49").Value = "  +1.91%  "

# Row 50
$ws.Range("E50").Value = "  +2.11%  "

# Row 51
$ws.Range("D51").Value = "79.45"
$ws.Range("E51").Value = "  +3.32%  "

# Restore default (Normal) style on the data range so only displayed values change
$ws.Range("D2:E51").Style = "Normal"

Write-Host "Cryptos list updated"